# Corrected error in 2019 figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 20: add the missing C20 = D20-B20 formula ---
$ws.Range("C20").Formula = "=D20-B20"

# --- Row 22: fix the D22 figure and add the C22 = D22-B22 formula ---
# (E22 = -D22 will recalculate automatically once D22 changes)
$ws.Range("D22").Value = 54736
$ws.Range("C22").Formula = "=D22-B22"

# --- View: move the visible top-left window and selected range ---
$ws.Range("B15:J25").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 5

Write-Output "done"
